$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 294, pushing the existing rows 294-347 down to 296-349.
$ws.Rows.Item(294).Insert()
$ws.Rows.Item(294).Insert()

# --- New row 294: Naranja / Valencia / Primera ---
$ws.Range("A294").Value2 = 7
$ws.Range("B294").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C294").Value2 = "Ñuble"
$ws.Range("D294").Value2 = 44641
$ws.Range("E294").Value2 = 16
$ws.Range("F294").Value2 = "Fruta"
$ws.Range("G294").Value2 = 100102
$ws.Range("H294").Value2 = "Cítricos"
$ws.Range("I294").Value2 = 100102005
$ws.Range("J294").Value2 = "Naranja"
$ws.Range("K294").Value2 = "Valencia"
$ws.Range("L294").Value2 = "Primera"
$ws.Range("M294").Value2 = 120
$ws.Range("N294").Value2 = 10000
$ws.Range("O294").Value2 = 11000
$ws.Range("P294").Value2 = 10500
$ws.Range("Q294").Value2 = "$/bandeja 15 kilos granel"
$ws.Range("R294").Value2 = "Región de O'Higgins"
$ws.Range("S294").Value2 = 700
$ws.Range("T294").Value2 = 15

# --- New row 295: Naranja / Valencia / Segunda ---
$ws.Range("A295").Value2 = 7
$ws.Range("B295").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C295").Value2 = "Ñuble"
$ws.Range("D295").Value2 = 44641
$ws.Range("E295").Value2 = 16
$ws.Range("F295").Value2 = "Fruta"
$ws.Range("G295").Value2 = 100102
$ws.Range("H295").Value2 = "Cítricos"
$ws.Range("I295").Value2 = 100102005
$ws.Range("J295").Value2 = "Naranja"
$ws.Range("K295").Value2 = "Valencia"
$ws.Range("L295").Value2 = "Segunda"
$ws.Range("M295").Value2 = 60
$ws.Range("N295").Value2 = 9000
$ws.Range("O295").Value2 = 9000
$ws.Range("P295").Value2 = 9000
$ws.Range("Q295").Value2 = "$/bandeja 15 kilos granel"
$ws.Range("R295").Value2 = "Región de O'Higgins"
$ws.Range("S295").Value2 = 600
$ws.Range("T295").Value2 = 15
